$d = $word.ActiveDocument

function Merge-Text([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Find.Execute failed for: $oldText"
    }
}

# ---------------------------------------------------------------------------
# 1) Merge runs that were split mid-sentence back into a single run (no
#    visible text change - these are no-op replacements that coalesce the
#    underlying <w:r> elements).
# ---------------------------------------------------------------------------

Merge-Text "Das System muss die für eine Wahl relevanten statischen Daten (Wahlbezirke, Wahlkreise, Kandidaten, Landeslisten, Parteien, Bundesländer) und dynamischen Daten (Stimmanzahl) für mehrere Wahlperioden (mindestens 2) speichern können, für zurückliegende Wahlen erfolgt die Speicherung der Stimmen in aggregierter Form auf Wahlkreisebene." `
            "Das System muss die für eine Wahl relevanten statischen Daten (Wahlbezirke, Wahlkreise, Kandidaten, Landeslisten, Parteien, Bundesländer) und dynamischen Daten (Stimmanzahl) für mehrere Wahlperioden (mindestens 2) speichern können, für zurückliegende Wahlen erfolgt die Speicherung der Stimmen in aggregierter Form auf Wahlkreisebene."

Merge-Text "Bei der Speicherung der Wahldaten muss sichergestellt sein, dass Daten weder verloren gehen, noch manipuliert oder unberechtigt ausgelesen werden können." `
            "Bei der Speicherung der Wahldaten muss sichergestellt sein, dass Daten weder verloren gehen, noch manipuliert oder unberechtigt ausgelesen werden können."

Merge-Text "Das System muss die korrekte Berechnung und tabellarische Anzeige der gewählten Kandidaten für den Bundestag auf Basis der abgegebenen Stimmen ermöglichen. Diese Berechnung soll auch auf Basis unvollständiger Daten (Hochrechnung) möglich sein; dafür gespeicherte Zwischenergebnisse auf Wahlkreisebene sollen durch den Benutzer angesteuert oder nach einem automatischen Verfahren neu berechnet werden." `
            "Das System muss die korrekte Berechnung und tabellarische Anzeige der gewählten Kandidaten für den Bundestag auf Basis der abgegebenen Stimmen ermöglichen. Diese Berechnung soll auch auf Basis unvollständiger Daten (Hochrechnung) möglich sein; dafür gespeicherte Zwischenergebnisse auf Wahlkreisebene sollen durch den Benutzer angesteuert oder nach einem automatischen Verfahren neu berechnet werden."

Merge-Text "Das System muss für Erststimmen auf Wahlkreisebene und für Zweitstimmen auf Landesebene die prozentuale Stimmverteilung berechnen und in tabellarischer und graphischer Form ausgeben können." `
            "Das System muss für Erststimmen auf Wahlkreisebene und für Zweitstimmen auf Landesebene die prozentuale Stimmverteilung berechnen und in tabellarischer und graphischer Form ausgeben können."

Merge-Text "Das System soll verschiedene Analysemöglichkeiten und die Ausgabe der Ergebnisse derselben in tabellarischer und graphischer Form bieten:" `
            "Das System soll verschiedene Analysemöglichkeiten und die Ausgabe der Ergebnisse derselben in tabellarischer und graphischer Form bieten:"

Merge-Text "Es soll ein Vergleich der Wahlbeteiligung zwischen Wahlkreisen/Bundesländern sowie über mehrere Wahlen hinweg möglich sein." `
            "Es soll ein Vergleich der Wahlbeteiligung zwischen Wahlkreisen/Bundesländern sowie über mehrere Wahlen hinweg möglich sein."

Merge-Text "Das System kann eine Möglichkeit zur Pflege der für die Wahl nötigen statischen Informationen (Wahlkreise, Wahlbezirke, Kandidaten, Landeslisten) bieten." `
            "Das System kann eine Möglichkeit zur Pflege der für die Wahl nötigen statischen Informationen (Wahlkreise, Wahlbezirke, Kandidaten, Landeslisten) bieten."

Merge-Text "Die technische Umsetzung des Systems erfolgt durch eine Browser-basierte Client-Server-Lösung." `
            "Die technische Umsetzung des Systems erfolgt durch eine Browser-basierte Client-Server-Lösung."

Merge-Text "Der Serverteil besteht aus einer relationalen Datenbank sowie einem auf Java basierenden Backend, das für die Datenpflege und Generierung der Ergebnisse und Analysen verantwortlich ist und durch ein Webframework die Browser-basierte Oberfläche anbietet. Die Berechnung von Ergebnissen und Analysen erfolgt dabei soweit wie möglich durch SQL-Statements und unter Ausnutzung der durch das gewählte Datenbanksystem angebotenen Optimierungsmöglichkeiten." `
            "Der Serverteil besteht aus einer relationalen Datenbank sowie einem auf Java basierenden Backend, das für die Datenpflege und Generierung der Ergebnisse und Analysen verantwortlich ist und durch ein Webframework die Browser-basierte Oberfläche anbietet. Die Berechnung von Ergebnissen und Analysen erfolgt dabei soweit wie möglich durch SQL-Statements und unter Ausnutzung der durch das gewählte Datenbanksystem angebotenen Optimierungsmöglichkeiten."

Merge-Text "Bei Design und Umsetzung der Architektur wird die Möglichkeit der Parallelisierung der Lösung im Sinne von Hochverfügbarkeit und gleichzeitigen Massenzugriffen berücksichtigt und soweit möglich vorbereitet." `
            "Bei Design und Umsetzung der Architektur wird die Möglichkeit der Parallelisierung der Lösung im Sinne von Hochverfügbarkeit und gleichzeitigen Massenzugriffen berücksichtigt und soweit möglich vorbereitet."

Merge-Text " zu garantieren und den Nutzern die Möglichkeit geben alle Ergebnisse auf einen Blick erfassen zu können." `
            " zu garantieren und den Nutzern die Möglichkeit geben alle Ergebnisse auf einen Blick erfassen zu können."

Merge-Text "Ein besonderes Konzept für die Zugangsberechtigung wird nicht erstellt, aggregierte Daten ab Wahlkreisebene sind für alle Benutzer ohne Login oder sonstige Authentifizierung zugänglich, Einzelstimmen können nicht eingesehen werden." `
            "Ein besonderes Konzept für die Zugangsberechtigung wird nicht erstellt, aggregierte Daten ab Wahlkreisebene sind für alle Benutzer ohne Login oder sonstige Authentifizierung zugänglich, Einzelstimmen können nicht eingesehen werden."

Merge-Text "Das System darf auf keinen Fall die Zuordnung von Wählern zu Stimmzetteln ermöglichen." `
            "Das System darf auf keinen Fall die Zuordnung von Wählern zu Stimmzetteln ermöglichen."

Merge-Text "Das Wahlinformationssystem dient ausschließlich Informationszwecken, es soll von den Wahlvorständen nicht als Ersatz für die Auszählung der Stimmen genutzt werden." `
            "Das Wahlinformationssystem dient ausschließlich Informationszwecken, es soll von den Wahlvorständen nicht als Ersatz für die Auszählung der Stimmen genutzt werden."

Merge-Text "Das System dient nicht als Plattform für die (Selbst-)Darstellung von Kandidaten oder Parteien." `
            "Das System dient nicht als Plattform für die (Selbst-)Darstellung von Kandidaten oder Parteien."

# ---------------------------------------------------------------------------
# 2) Actual text changes
# ---------------------------------------------------------------------------

Merge-Text "Das System muss eine Webservice-Schnittstelle zum Einpflegen von Stimmzetteln anbieten." `
            "Das System muss die Abgabe von Stimmzetteln über eine Web-Oberfläche ermöglichen."

Merge-Text "Als zusätzlicher Zugriffsweg wird eine REST-basierte Webservice-Schnittstelle implementiert, durch die abgegebene Stimmzettel eingepflegt werden können. Für die Massenladung von Stimmzetteln wird ein zusätzliches Client-Programm auf Java-Basis entwickelt, das die angelieferten Stimmzettel per Webservice in das System einpflegt." `
            "Für die Massenladung von Stimmzetteln wird ein Client-Programm auf Java-Basis entwickelt, das die digital angelieferten Stimmzettel in das System einpflegt."

Merge-Text "Das System soll es einem Wähler nicht ermöglichen, seine Stimme direkt über die Weboberfläche abzugeben, für die Abgabe von Stimmen ist weiterhin das Wahllokal zuständig." `
            "Das System soll es einem Wähler ermöglichen, seine Stimme über eine im Wahllokal bereitgestellte Weboberfläche elektronisch abzugeben."

# ---------------------------------------------------------------------------
# 3) Remove the whole "Glossar" section (heading + 2 glossary entries + the
#    trailing empty paragraph before the _GoBack bookmark paragraph).
# ---------------------------------------------------------------------------

$rng = $d.Content
$startFound = $rng.Find.Execute("Glossar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $startFound) { throw "Could not find Glossar heading" }
$glossarStart = $rng.Start

$rng2 = $d.Content
$endFound = $rng2.Find.Execute("tiert werden.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $endFound) { throw "Could not find end of REST glossary entry" }
$afterRest = $rng2.End

# Also remove the trailing empty "ind left=0" paragraph that sits between the
# REST entry and the bookmark paragraph - it is two paragraph marks further on
# (end of the REST paragraph itself, then the whole empty paragraph).
$glossarEnd = $afterRest + 2

$delRange = $d.Range($glossarStart, $glossarEnd)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 4) Remove the old `_GoBack` bookmark (was sitting in its own paragraph near
#    the former Glossar section - that paragraph now becomes a plain empty
#    "Listenabsatz" paragraph).
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 5) Add the new `_GoBack` bookmark around "Zielsetzung" (Word stamps this at
#    the location of the most recent edit).
# ---------------------------------------------------------------------------

$zielRng = $d.Content
$zFound = $zielRng.Find.Execute("Zielsetzung", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $zFound) { throw "Could not find Zielsetzung" }
$d.Bookmarks.Add("_GoBack", $zielRng)

# ---------------------------------------------------------------------------
# 6) Update the cached field results in the header (page count + date).
# ---------------------------------------------------------------------------

$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRng = $hdr.Range

$pageFound = $hdrRng.Find.Execute("Seite " + [char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$hdrRng2 = $hdr.Range
$found1 = $hdrRng2.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2)

$hdrRng3 = $hdr.Range
$found2 = $hdrRng3.Find.Execute("02.11.2012", $true, $false, $false, $false, $false, $true, 1, $false, "18.01.2013", 2)
